$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the Volume/Number and report date-range headers ---
$ws.Range("A8").Value = "Volume 31   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/16/2024  Through  12/22/2024"

# --- 2) Insert a new blank row before row 56 (pushes the trailing two rows down by one) ---
$ws.Rows("56:56").Insert()

# --- 3) Cells that change numeric <-> placeholder-text type need their number format changed too.
#        Donor cells with the right format are pasted in (format-only) before/after setting the value
#        so the stored style index matches a plain numeric / plain General-text cell respectively. ---
$ws.Range("G14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = 1

$ws.Range("G14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value = 1

$ws.Range("G14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 4

$ws.Range("E16").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = 0

$ws.Range("A14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("G14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 2

$excel.CutCopyMode = $false

# --- 4) Plain numeric value updates (no style/type change) ---
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 20
$ws.Range("K14").Value = 185.714285714286
$ws.Range("L14").Value = -4.761904761904
$ws.Range("M14").Value = -28.571428571428
$ws.Range("N14").Value = -72.972972972973

# Row 15
$ws.Range("L15").Value = 2.777777777777
$ws.Range("N15").Value = -56.976744186046

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 25
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = 8.695652173913
$ws.Range("I16").Value = 254
$ws.Range("J16").Value = 271
$ws.Range("K16").Value = -6.273062730627
$ws.Range("L16").Value = -25.730994152046
$ws.Range("M16").Value = -48.057259713701
$ws.Range("N16").Value = -88.449295134151

# Row 17
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 42
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = -2.325581395348
$ws.Range("I17").Value = 630
$ws.Range("J17").Value = 671
$ws.Range("K17").Value = -6.110283159463
$ws.Range("L17").Value = -10.638297872340
$ws.Range("M17").Value = 11.900532859680
$ws.Range("N17").Value = -43.497757847533

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 139
$ws.Range("J18").Value = 125
$ws.Range("K18").Value = 11.2
$ws.Range("L18").Value = -27.225130890052
$ws.Range("M18").Value = -49.637681159420
$ws.Range("N18").Value = -82.018111254851

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -28.125
$ws.Range("I19").Value = 346
$ws.Range("J19").Value = 337
$ws.Range("K19").Value = 2.670623145400
$ws.Range("L19").Value = -14.567901234567
$ws.Range("M19").Value = -21.896162528216
$ws.Range("N19").Value = -62.994652406417

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -44.444444444444
$ws.Range("I20").Value = 138
$ws.Range("J20").Value = 148
$ws.Range("K20").Value = -6.756756756756
$ws.Range("L20").Value = -25
$ws.Range("M20").Value = -3.496503496503
$ws.Range("N20").Value = -80.285714285714

# Row 21
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -12.903225806451
$ws.Range("G21").Value = 121
$ws.Range("H21").Value = -12.396694214876
$ws.Range("I21").Value = 1564
$ws.Range("J21").Value = 1595
$ws.Range("K21").Value = -1.943573667711
$ws.Range("L21").Value = -16.985138004246
$ws.Range("M21").Value = -20.850202429149
$ws.Range("N21").Value = -73.410404624277

# Row 22
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 53
$ws.Range("K22").Value = 13.207547169811
$ws.Range("L22").Value = 22.448979591836
$ws.Range("M22").Value = 1.694915254237

# Row 23
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 66.666666666666
$ws.Range("G23").Value = 24
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 337
$ws.Range("J23").Value = 338
$ws.Range("K23").Value = -0.295857988165
$ws.Range("L23").Value = -10.372340425531
$ws.Range("M23").Value = 35.887096774193

# Row 24
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 24
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = -2.777777777777
$ws.Range("I24").Value = 1172
$ws.Range("J24").Value = 1071
$ws.Range("K24").Value = 9.430438842203
$ws.Range("L24").Value = 3.808680248007
$ws.Range("M24").Value = 19.591836734693

# Row 25
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -44.444444444444
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = -48.936170212766
$ws.Range("I25").Value = 307
$ws.Range("J25").Value = 298
$ws.Range("K25").Value = 3.020134228187
$ws.Range("L25").Value = -10.233918128655

# Row 26
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = 11.764705882352
$ws.Range("F26").Value = 84
$ws.Range("G26").Value = 55
$ws.Range("H26").Value = 52.727272727272
$ws.Range("I26").Value = 894
$ws.Range("J26").Value = 928
$ws.Range("K26").Value = -3.663793103448
$ws.Range("L26").Value = -3.455723542116
$ws.Range("M26").Value = -29.606299212598

# Row 27
$ws.Range("L27").Value = -18.644067796610

# Row 28
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -62.5
$ws.Range("I28").Value = 71
$ws.Range("J28").Value = 81
$ws.Range("K28").Value = -12.345679012345
$ws.Range("L28").Value = -17.441860465116

# Row 29
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 50
$ws.Range("I29").Value = 59
$ws.Range("J29").Value = 47
$ws.Range("K29").Value = 25.531914893617
$ws.Range("L29").Value = -22.368421052631
$ws.Range("M29").Value = -43.269230769230
$ws.Range("N29").Value = -81.901840490797

# Row 30
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 25
$ws.Range("I30").Value = 49
$ws.Range("J30").Value = 43
$ws.Range("K30").Value = 13.953488372093
$ws.Range("L30").Value = -18.333333333333
$ws.Range("M30").Value = -41.666666666666
$ws.Range("N30").Value = -83.501683501683
